$d = $word.ActiveDocument

# --- Locate the target paragraph: the literature-review body paragraph that
# --- follows the "Theory and Literature Review" heading (not the earlier
# --- "Background, Aims and Objectives" instructions paragraph which has the
# --- same opening sentence).
$headingRng = $d.Content
$headingRng.Find.Execute("Theory and Literature Review", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$searchRng = $d.Range($headingRng.End, $d.Content.End)
$searchRng.Find.Execute("This section should be a thorough examination", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bodyPara = $searchRng.Paragraphs(1)
$startPos = $bodyPara.Range.Start

# The paragraph is immediately followed by three empty paragraphs that are
# also being folded into the new content block (the new text replaces all
# four original paragraphs).
$spanRange = $d.Range($startPos, $d.Content.End)
$fourthPara = $spanRange.Paragraphs(4)
$endPos = $fourthPara.Range.End

$targetRange = $d.Range($startPos, $endPos)

$newXml = @'
<w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">The exponential growth in computational demands, driven by applications in machine learning, multimedia processing, and big data analytics, has strained traditional digital design paradigms. Classical computing architectures prioritise precision and exactness, which come at the cost of increased power consumption, area usage, and latency. With the diminishing benefits of Moore’s Law and the rising need for energy-efficient hardware, approximate computing has emerged as a transformative approach to hardware design. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>Approximate computing operates on the principle that not all applications require perfect accuracy. Many domains, especially those involving human perception or probabilistic outcomes, can tolerate small errors without significant degradation in performance. By trade-off of accuracy, approximate computing reduces hardware complexity, resulting in substantial improvements in energy efficiency, and processing speed.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">At the heart of this model shift are arithmetic units like adders and multipliers which constitute a significant portion of computational workloads in digital systems. Optimising these units for approximate computing forms the core of this paper’s contributions. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Approximate Adders</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">Adders are a </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">fundamental component in </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">digital </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>circuits, responsible for executing arithmetic operations that often dominate computational workload. Traditional adder designs prioritise accuracy, however, approximate adders introduce intentional inaccuracies to achieve resource savings</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>A proposed approximation approach using Lower-Part OR-based Approximate Adders [1] aligns with similar research</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>introducing the concept of approximate adders as a means to trade</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">off accuracy for reduced power consumptions and area in energy-efficient VLSI systems. Ramasamy et al. proposed a carry-based approximate full adder, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t xml:space="preserve">demonstrating that bypassing the carry propagation chain in the least significant bits (LSB) can drastically improve speed and reduce area at the cost of negligible error [reference]. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Approximate Multipliers</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
        <w:t>Multiplication is a computationally intensive operation, making approximate multipliers a critical focus for energy-efficient design. Approximate multipliers reduce the complexity of partial product summation, which directly impacts delay and power consumption.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Approximate Matrix Multiplication</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Compressor-Based Approximate Multipliers</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Novel Compressor Designs</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Once-Through Multiplier Architecture CAM2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>FPGA-Based Approximate Multipliers</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="2268"/>
          <w:tab w:val="right" w:pos="10093"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Source Sans Pro"/>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Trade-Offs and Error Analysis</w:t>
      </w:r>
    </w:p>

'@

$payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($payload)

Write-Output "Literature review section expanded."
